$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-10-09 Monday" "2023-10-10 Tuesday"

Replace-Text "29÷4=7, 1" "47÷8=5, 7"
Replace-Text "82÷7=11, 5" "39÷8=4, 7"
Replace-Text "30÷5=6, 0" "16÷9=1, 7"
Replace-Text "74÷6=12, 2" "88÷3=29, 1"
Replace-Text "38÷9=4, 2" "68÷4=17, 0"

Replace-Text "23÷2=11, 1" "94÷7=13, 3"
Replace-Text "66÷9=7, 3" "56÷2=28, 0"
Replace-Text "15÷6=2, 3" "98÷3=32, 2"
Replace-Text "87÷9=9, 6" "35÷4=8, 3"
Replace-Text "19÷3=6, 1" "42÷7=6, 0"

Replace-Text "35÷9=3, 8" "89÷7=12, 5"
Replace-Text "27÷5=5, 2" "92÷3=30, 2"
Replace-Text "27÷3=9, 0" "95÷7=13, 4"
Replace-Text "64÷8=8, 0" "81÷8=10, 1"
Replace-Text "76÷8=9, 4" "17÷9=1, 8"

Replace-Text "77÷6=12, 5" "68÷8=8, 4"
Replace-Text "59÷8=7, 3" "89÷8=11, 1"
Replace-Text "83÷8=10, 3" "60÷3=20, 0"
Replace-Text "69÷8=8, 5" "58÷7=8, 2"
Replace-Text "67÷2=33, 1" "73÷3=24, 1"

Replace-Text "65÷8=8, 1" "73÷7=10, 3"
Replace-Text "61÷4=15, 1" "97÷7=13, 6"
Replace-Text "91÷7=13, 0" "73÷5=14, 3"
Replace-Text "24÷4=6, 0" "87÷6=14, 3"
Replace-Text "13÷7=1, 6" "70÷3=23, 1"

Write-Output "done"
